$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook is a "handback status" report with three sheets:
#   Overview (sheet1), zh-cn (sheet2), de-de (sheet3)
# A CI re-run regenerated the report against a new pair of source/target
# files (UUID-named .md source docs), so every reference to the old UUIDs
# and the old generation timestamps needs to move to the new ones.
#
#   old "00ea13d2-1419-4f8c-9d94-be771667e003" -> new "f8f155ac-b6d9-4825-800a-3cac39ded6bb"
#   old "51bcdee1-b9d0-4a67-a8b4-136896f59d72" -> new "ffffeb5154df-c573-43cc-9662-0c4f1ab57e07"
#   old xlf hash "064b06a2cb5c5f99c40fbdc681e8a627cbce8380" (file 1) -> new "8fe65fe2b16217659ca8b752011a1fbee6de9245"
#   old xlf hash "989e6333fa360a38a15ecab03460da3f55733caf" (file 2) -> new "8fe65fe2b16217659ca8b752011a1fbee6de9245" (same as file 1 this run)
#
# Net effect: both rows of the zh-cn/de-de sheets now reference the SAME
# handback xlf artifact (the CI run only produced one xlf per language this
# time), and all the "generate/handback" timestamps move forward to the new
# run's timestamps.
# ---------------------------------------------------------------------------

$oldUuid1 = "00ea13d2-1419-4f8c-9d94-be771667e003"
$newUuid1 = "f8f155ac-b6d9-4825-800a-3cac39ded6bb"
$oldUuid2 = "51bcdee1-b9d0-4a67-a8b4-136896f59d72"
$newUuid2 = "ffffeb5154df-c573-43cc-9662-0c4f1ab57e07"

$newMd1 = $newUuid1 + ".md"
$newMd2 = $newUuid2 + ".md"
$newMdPath1 = "e2e\" + $newMd1
$newMdPath2 = "e2e\" + $newMd2

$newXlfZh = $newUuid1 + ".8fe65fe2b16217659ca8b752011a1fbee6de9245.zh-cn.xlf"
$newXlfDe = $newUuid1 + ".8fe65fe2b16217659ca8b752011a1fbee6de9245.de-de.xlf"

$newGenerateDate = "2016-08-16 23:01:30"
$newZhGenDate     = "2016-08-16 23:01:25"
$newZhHandbackDate= "2016-08-16 23:01:41"
$newDeHandbackDate= "2016-08-16 23:01:48"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newMdPath1
$wsOverview.Range("G2").Value = $newGenerateDate

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newMdPath2
$wsOverview.Range("G3").Value = $newGenerateDate

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') {
        $h.TextToDisplay = $newMdPath1
    } elseif ($h.Range.Address() -eq '$B$3') {
        $h.TextToDisplay = $newMdPath2
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newZhGenDate
$wsZh.Range("I2").Value = $newMd1
$wsZh.Range("J2").Value = $newXlfZh
$wsZh.Range("K2").Value = $newZhHandbackDate

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newZhGenDate
$wsZh.Range("I3").Value = $newMd2
$wsZh.Range("J3").Value = $newXlfZh
$wsZh.Range("K3").Value = $newZhHandbackDate

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = $newMd1
    } elseif ($h.Range.Address() -eq '$I$2') {
        $h.TextToDisplay = $newMd1
    } elseif ($h.Range.Address() -eq '$A$3') {
        $h.TextToDisplay = $newMd2
    } elseif ($h.Range.Address() -eq '$I$3') {
        $h.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newGenerateDate
$wsDe.Range("I2").Value = $newMd1
$wsDe.Range("J2").Value = $newXlfDe
$wsDe.Range("K2").Value = $newDeHandbackDate

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newGenerateDate
$wsDe.Range("I3").Value = $newMd2
$wsDe.Range("J3").Value = $newXlfDe
$wsDe.Range("K3").Value = $newDeHandbackDate

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = $newMd1
    } elseif ($h.Range.Address() -eq '$I$2') {
        $h.TextToDisplay = $newMd1
    } elseif ($h.Range.Address() -eq '$A$3') {
        $h.TextToDisplay = $newMd2
    } elseif ($h.Range.Address() -eq '$I$3') {
        $h.TextToDisplay = $newMd2
    }
}
